$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.838.61'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '2.439.47'
$ws.Range('E3').Value = '  -1.30%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '558.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '162.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.37%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  -0.61%  '
$ws.Range('E9').Value = '  +7.96%  '
$ws.Range('E10').Value = '  -2.18%  '
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.60'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000177'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.87%  '
$ws.Range('D14').Value = '68.724.20'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = '2.887.96'
$ws.Range('E15').Value = '  -0.61%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.25'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.53%  '
$ws.Range('D17').Value = '2.440.87'
$ws.Range('E17').Value = '  -1.00%  '
$ws.Range('E18').Value = '  -0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '339.18'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.30%  '
$ws.Range('E20').Value = '  +0.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.83'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.97%  '
$ws.Range('E22').Value = '  +3.01%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').Value = '2.567.19'
$ws.Range('E26').Value = '  -1.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.963'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.72%  '
$ws.Range('D29').Value = '0.0₃0820'
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.45%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('E32').Value = '  +1.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '429.72'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('E34').Value = '  -1.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '160.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.01'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('E39').Value = '  -2.35%  '
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('E41').Value = '  +2.08%  '
$ws.Range('E42').Value = '  -2.34%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.06'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.36%  '
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.32'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '130.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('E47').Value = '  +0.12%  '
$ws.Range('E48').Value = '  -1.08%  '
$ws.Range('E49').Value = '  -1.02%  '
$ws.Range('E50').Value = '  +2.78%  '
$ws.Range('E51').Value = '  +0.21%  '
